$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 19-22: the B column currently holds "1" and C column holds the real value.
# Move the C column value into B, and clear the C column (so B ends up holding
# the value that used to be in C, and C becomes empty).
for ($r = 19; $r -le 22; $r++) {
    $cValue = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($r, 2).Value = $cValue
    $ws.Cells.Item($r, 3).Value = $null
}

# Update the selected/active cell to C20 as reflected in the sheet view.
$ws.Range("C20").Select()
